$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 0
    3  = 2
    4  = 4
    5  = 3
    6  = 2
    7  = 4
    8  = 2
    9  = 3
    10 = 0
    11 = 3
    12 = 5
    13 = 2
    14 = 1
    15 = 0
    16 = 1
    17 = 2
    18 = 2
    19 = 4
    20 = 8
    21 = 1
    22 = 1
    23 = 1
    24 = 2
    25 = 1
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 7).Value = $values[$row]
}
